# change scheduling and correct get job type in linkedin module
# Update the last existing row (A58) and append the newly scraped job links.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A58").Value = "https://cryptocurrencyjobs.co/operations/futuremoney-group-operation-events-coordinator/"
$ws.Range("A59").Value = "https://cryptocurrencyjobs.co/engineering/chainlink-incident-responder-apac/"
$ws.Range("A60").Value = "https://cryptocurrencyjobs.co/marketing/chainlink-market-strategy-manager-capital-markets/"
$ws.Range("A61").Value = "https://cryptocurrencyjobs.co/engineering/chainlink-principal-engineer/"
$ws.Range("A62").Value = "https://cryptocurrencyjobs.co/engineering/chainlink-research-engineer/"
$ws.Range("A63").Value = "https://cryptocurrencyjobs.co/sales/chainlink-senior-data-partnerships-manager/"
$ws.Range("A64").Value = "https://cryptocurrencyjobs.co/engineering/chainlink-senior-site-reliability-engineer/"
$ws.Range("A65").Value = "https://cryptocurrencyjobs.co/operations/chainlink-strategic-finance-manager/"
$ws.Range("A66").Value = "https://cryptocurrencyjobs.co/engineering/chainlink-engineering-manager-blockchain-integrations/"
$ws.Range("A67").Value = "https://cryptocurrencyjobs.co/sales/chainlink-technical-account-manager/"
$ws.Range("A68").Value = "https://cryptocurrencyjobs.co/engineering/ontropy-senior-blockchain-engineer/"
$ws.Range("A69").Value = "https://cryptocurrencyjobs.co/sales/gelato-network-business-development-internship/"
$ws.Range("A70").Value = "https://cryptocurrencyjobs.co/operations/chainlink-senior-people-business-partner/"
$ws.Range("A71").Value = "https://cryptocurrencyjobs.co/marketing/gridplus-marketing-director/"
$ws.Range("A72").Value = "https://cryptocurrencyjobs.co/marketing/blockdaemon-growth-strategy-lead/"
